$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.431.60"
$ws.Range("E2").Value = "  +4.49%  "
$ws.Range("D3").Value = "3.131.06"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'243.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.16%  "
$ws.Range("D6").Value = "'611.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("E7").Value = "  +1.93%  "
$ws.Range("D8").Value = "'0.383"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.99%  "
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").Value = "3.126.74"
$ws.Range("E10").Value = "  +0.09%  "
$ws.Range("D11").Value = "'0.783"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.51%  "
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").Value = "97.119.38"
$ws.Range("D14").Value = "'0.0000240"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.32%  "
$ws.Range("D15").Value = "'5.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.15%  "
$ws.Range("D16").Value = "'34.00"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.90%  "
$ws.Range("D17").Value = "3.714.78"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").Value = "3.134.22"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").Value = "'3.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.87%  "
$ws.Range("D20").Value = "'513.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +16.04%  "
$ws.Range("D21").Value = "'14.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.02%  "
$ws.Range("D22").Value = "'5.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.37%  "
$ws.Range("D23").Value = "'0.0000193"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.62%  "
$ws.Range("D24").Value = "'8.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.34%  "
$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D25").Value = "'5.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.05%  "
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").Value = "'88.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.17%  "
$ws.Range("D27").Value = "'11.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -10.78%  "
$ws.Range("D28").Value = "3.303.52"
$ws.Range("E28").Value = "  +0.31%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  +3.12%  "
$ws.Range("E31").Value = "  -3.42%  "
$ws.Range("D32").Value = "'0.124"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.57%  "
$ws.Range("D33").Value = "'0.972"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.57%  "
$ws.Range("D34").Value = "'8.99"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.98%  "
$ws.Range("D35").Value = "'26.57"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.08%  "
$ws.Range("D36").Value = "'0.152"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.42%  "
$ws.Range("D37").Value = "'7.32"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -9.48%  "
$ws.Range("E38").Value = "  -1.85%  "
$ws.Range("D39").Value = "'24.22"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.94%  "
$ws.Range("D40").Value = "'470.15"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.08%  "
$ws.Range("D41").Value = "'0.436"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.98%  "
$ws.Range("D42").Value = "'3.60"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -9.42%  "
$ws.Range("D43").Value = "'1.22"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.15%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").Value = "'3.13"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.91%  "
$ws.Range("D46").Value = "'161.94"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.06%  "
$ws.Range("D47").Value = "'1.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.40%  "
$ws.Range("D48").Value = "'0.698"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.45%  "
$ws.Range("D49").Value = "'4.48"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("D50").Value = "'44.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.39%  "
$ws.Range("D51").Value = "'0.998"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.00%  "
